$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 15611
$ws.Cells.Item(2, 2).Value = "dan"
$ws.Cells.Item(2, 3).ClearContents()
$ws.Cells.Item(2, 4).Value = 87
$ws.Cells.Item(2, 5).ClearContents()
$ws.Cells.Item(2, 6).Value = 74
$ws.Cells.Item(2, 7).ClearContents()
$ws.Cells.Item(2, 8).ClearContents()
$ws.Cells.Item(2, 9).ClearContents()
$ws.Cells.Item(2, 10).Value = 76
$ws.Cells.Item(2, 11).Value = 237
$ws.Cells.Item(2, 12).Value = 29.625
$ws.Cells.Item(2, 13).Value = "E"

# Row 3
$ws.Cells.Item(3, 1).Value = 15612
$ws.Cells.Item(3, 2).Value = "tes"
$ws.Cells.Item(3, 3).ClearContents()
$ws.Cells.Item(3, 4).Value = 89
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(3, 6).Value = 48
$ws.Cells.Item(3, 7).ClearContents()
$ws.Cells.Item(3, 8).ClearContents()
$ws.Cells.Item(3, 9).ClearContents()
$ws.Cells.Item(3, 10).Value = 78
$ws.Cells.Item(3, 11).Value = 215
$ws.Cells.Item(3, 12).Value = 26.875
$ws.Cells.Item(3, 13).Value = "E"
